$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.732.74"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").Value = "2.095.07"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.5158"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("D8").Value = "'0.4379"
$ws.Range("E8").Value = "  -3.55%  "
$ws.Range("D9").Value = "'53.02"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").Value = "'0.09216"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").Value = "'24.80"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "2.105.32"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").Value = "'8.241"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").Value = "'6.761"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "'0.00001150"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "'20.77"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "'0.06662"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "'6.205"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("D23").Value = "29.765.36"
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").Value = "'12.47"
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("D25").Value = "'2.319"
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D26").Value = "2.350.10"
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("D27").Value = "'21.93"
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").Value = "'2.513"
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("D29").Value = "'161.16"
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").Value = "'132.96"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("E31").Value = "  -6.35%  "
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").Value = "'1.626"
$ws.Range("E33").Value = "  -4.56%  "
$ws.Range("D34").Value = "'6.167"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("D36").Value = "'6.291"
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("D37").Value = "'10.19"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").Value = "'0.02574"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("D39").Value = "'0.7108"
$ws.Range("E39").Value = "  +2.21%  "
$ws.Range("D40").Value = "'0.06721"
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.323"
$ws.Range("E42").Value = "  +3.96%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.2225"
$ws.Range("E43").Value = "  -5.27%  "
$ws.Range("D44").Value = "'0.7015"
$ws.Range("E44").Value = "  +8.23%  "
$ws.Range("D45").Value = "'14.29"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").Value = "'1.009"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "'2.314"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "'3.619"
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("D49").Value = "'0.00000000354"
$ws.Range("E49").Value = "  -5.30%  "
$ws.Range("D50").Value = "'1.219"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("D51").Value = "'81.94"
$ws.Range("E51").Value = "  -1.93%  "
